{"js": "const replacements = [\n  [\"2024-11-26 Tuesday\", \"2024-11-27 Wednesday\"],\n  [\"736\u00d74=2944\", \"191\u00d74=764\"],\n  [\"655\u00d78=5240\", \"996\u00d78=7968\"],\n  [\"261\u00d79=2349\", \"927\u00d77=6489\"],\n  [\"541\u00d74=2164\", \"431\u00d77=3017\"],\n  [\"733\u00d73=2199\", \"160\u00d74=640\"],\n  [\"500\u00d74=2000\", \"889\u00d72=1778\"],\n  [\"753\u00d77=5271\", \"641\u00d74=2564\"],\n  [\"205\u00d79=1845\", \"123\u00d74=492\"],\n  [\"772\u00d77=5404\", \"998\u00d76=5988\"],\n  [\"144\u00d73=432\", \"987\u00d73=2961\"],\n  [\"458\u00d74=1832\", \"162\u00d78=1296\"],\n  [\"890\u00d78=7120\", \"460\u00d77=3220\"],\n  [\"149\u00d75=745\", \"335\u00d78=2680\"],\n  [\"709\u00d76=4254\", \"615\u00d72=1230\"],\n  [\"971\u00d77=6797\", \"445\u00d78=3560\"],\n  [\"523\u00d75=2615\", \"399\u00d78=3192\"],\n  [\"420\u00d72=840\", \"994\u00d76=5964\"],\n  [\"200\u00d76=1200\", \"935\u00d74=3740\"],\n  [\"792\u00d74=3168\", \"720\u00d75=3600\"],\n  [\"371\u00d74=1484\", \"543\u00d73=1629\"],\n  [\"593\u00d74=2372\", \"360\u00d77=2520\"],\n  [\"860\u00d75=4300\", \"558\u00d73=1674\"],\n  [\"775\u00d74=3100\", \"409\u00d79=3681\"],\n  [\"918\u00d78=7344\", \"491\u00d77=3437\"],\n  [\"552\u00d75=2760\", \"853\u00d79=7677\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@('2024-11-26 Tuesday', '2024-11-27 Wednesday')\n    ,@('736\u00d74=2944', '191\u00d74=764')\n    ,@('655\u00d78=5240', '996\u00d78=7968')\n    ,@('261\u00d79=2349', '927\u00d77=6489')\n    ,@('541\u00d74=2164', '431\u00d77=3017')\n    ,@('733\u00d73=2199', '160\u00d74=640')\n    ,@('500\u00d74=2000', '889\u00d72=1778')\n    ,@('753\u00d77=5271', '641\u00d74=2564')\n    ,@('205\u00d79=1845', '123\u00d74=492')\n    ,@('772\u00d77=5404', '998\u00d76=5988')\n    ,@('144\u00d73=432', '987\u00d73=2961')\n    ,@('458\u00d74=1832', '162\u00d78=1296')\n    ,@('890\u00d78=7120', '460\u00d77=3220')\n    ,@('149\u00d75=745', '335\u00d78=2680')\n    ,@('709\u00d76=4254', '615\u00d72=1230')\n    ,@('971\u00d77=6797', '445\u00d78=3560')\n    ,@('523\u00d75=2615', '399\u00d78=3192')\n    ,@('420\u00d72=840', '994\u00d76=5964')\n    ,@('200\u00d76=1200', '935\u00d74=3740')\n    ,@('792\u00d74=3168', '720\u00d75=3600')\n    ,@('371\u00d74=1484', '543\u00d73=1629')\n    ,@('593\u00d74=2372', '360\u00d77=2520')\n    ,@('860\u00d75=4300', '558\u00d73=1674')\n    ,@('775\u00d74=3100', '409\u00d79=3681')\n    ,@('918\u00d78=7344', '491\u00d77=3437')\n    ,@('552\u00d75=2760', '853\u00d79=7677')\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $found = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        throw \"Text not found: $old\"\n    }\n}\n"}
